# The commit swaps the contents of ppt/theme/theme1.xml (used by the
# slide master -> "Integral" colours) and ppt/theme/theme2.xml (used by
# the notes master -> default "Office Theme" colours): after the edit,
# theme1.xml carries the plain "Office Theme" colour scheme and
# theme2.xml carries the "Integral" colour scheme.
#
# The PowerPoint object model only exposes the *active* (slide-master)
# theme's 12 colour slots for editing -- via Slide.ThemeColorScheme --
# so we drive that surface to repoint theme1.xml at the Office Theme
# palette.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function RGBInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Office Theme colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
$tcs.Colors(1).RGB  = RGBInt("000000")  # dk1
$tcs.Colors(2).RGB  = RGBInt("FFFFFF")  # lt1
$tcs.Colors(3).RGB  = RGBInt("44546A")  # dk2
$tcs.Colors(4).RGB  = RGBInt("E7E6E6")  # lt2
$tcs.Colors(5).RGB  = RGBInt("5B9BD5")  # accent1
$tcs.Colors(6).RGB  = RGBInt("ED7D31")  # accent2
$tcs.Colors(7).RGB  = RGBInt("A5A5A5")  # accent3
$tcs.Colors(8).RGB  = RGBInt("FFC000")  # accent4
$tcs.Colors(9).RGB  = RGBInt("4472C4")  # accent5
$tcs.Colors(10).RGB = RGBInt("70AD47")  # accent6
$tcs.Colors(11).RGB = RGBInt("0563C1")  # hlink
$tcs.Colors(12).RGB = RGBInt("954F72")  # folHlink
